# Remove the anchored logo picture ("Picture 2") that was embedded as the
# first run of the document's opening paragraph. The picture is a floating
# (anchored) drawing, so it shows up in the Shapes collection rather than
# InlineShapes; deleting the shape removes the whole <w:r><w:drawing>...
# run from the underlying XML while leaving the remaining runs (tab +
# company-name text) in that paragraph untouched.

$d = $word.ActiveDocument

for ($i = $d.Shapes.Count; $i -ge 1; $i--) {
    $shape = $d.Shapes.Item($i)
    if ($shape.Name -eq "Picture 2") {
        $shape.Delete()
    }
}
